$d = $word.ActiveDocument

# Locate the paragraph that ends the "2022-02-04:" note (the one that
# finishes with "...convert it into the DAGs."), so the new note can be
# appended right after it, before the pre-existing blank paragraph.
$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*convert it into the DAGs.*") {
        $targetIndex = $i
        break
    }
}

if ($targetIndex -eq -1) {
    throw "Could not find the '...convert it into the DAGs.' paragraph."
}

$targetParagraph = $d.Paragraphs.Item($targetIndex)

# Split off a brand-new (still empty) paragraph right after it; the
# document's existing blank paragraph gets pushed one further down.
$targetParagraph.Range.InsertParagraphAfter()

$newParagraph = $d.Paragraphs.Item($targetIndex + 1)

# Fill that freshly-minted paragraph with a blank paragraph followed by
# the new dated note, built from individual runs (mirroring how Word
# breaks typed text into runs / flags "hql" as a misspelling) via a
# WordprocessingML package fragment.
$xml = @"
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p/>
          <w:p>
            <w:r><w:t>2022-02-08:  Existing documentation on the ECL IR format can be found at the top of the C++ source code file</w:t></w:r>
            <w:r><w:t xml:space="preserve"> </w:t></w:r>
            <w:r><w:t>ecl/</w:t></w:r>
            <w:proofErr w:type="spellStart"/>
            <w:r><w:t>hql</w:t></w:r>
            <w:proofErr w:type="spellEnd"/>
            <w:r><w:t>/hqlir.cpp</w:t></w:r>
            <w:r><w:t>.</w:t></w:r>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
"@

$newParagraph.Range.InsertXML($xml) | Out-Null
